$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Player Info" sheet before "ODI Batting" ---
$battingSheetRef = $wb.Worksheets.Item("ODI Batting")
$infoSheet = $wb.Worksheets.Add($battingSheetRef)
$infoSheet.Name = "Player Info"

# Re-fetch fresh references by name (post-Add, to avoid any stale handles)
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 2. Populate "Player Info" sheet ---
$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# Write the numeric-looking ID as text (matches the source data convention
# where every cell - even numeric-looking ones - is stored as a string)
$infoSheet.Range("A2").NumberFormat = "@"
$infoSheet.Range("A2").Value = "4748"
$infoSheet.Range("A2").Style = "Normal"

$infoSheet.Range("B2").Value = "Odean Fabian Smith"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Right Arm Medium"

# Match the bold/bordered/centered header style used by the other sheets
# (copy formatting only, so the existing shared style gets reused instead
# of a near-duplicate style being minted)
$battingSheet.Range("A1").Copy()
$infoSheet.Range("A1:D1").PasteSpecial(-4122)

# --- 3. Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingMatchCodes = @{
    "D2" = "4519"
    "D3" = "4520"
    "D4" = "4522"
    "D5" = "4535"
    "D6" = "4536"
    "D7" = "4727"
    "D8" = "4731"
}
foreach ($addr in $battingMatchCodes.Keys) {
    $battingSheet.Range($addr).NumberFormat = "@"
    $battingSheet.Range($addr).Value = $battingMatchCodes[$addr]
    $battingSheet.Range($addr).Style = "Normal"
}

# --- 4. Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingMatchCodes = @{
    "B2" = "4519"
    "B3" = "4520"
    "B4" = "4522"
    "B5" = "4535"
    "B6" = "4536"
    "B7" = "4727"
    "B8" = "4731"
}
foreach ($addr in $bowlingMatchCodes.Keys) {
    $bowlingSheet.Range($addr).NumberFormat = "@"
    $bowlingSheet.Range($addr).Value = $bowlingMatchCodes[$addr]
    $bowlingSheet.Range($addr).Style = "Normal"
}
